$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds one (Primera/Segunda) pair of price rows per weekly
# observation, ordered from most recent to oldest starting at row 376.
# A new, more recent weekly observation is being added at the top of
# that block, pushing every existing observation down by two rows
# (376:377 -> 378:379, ..., 414:415 -> 416:417).

# Insert two blank rows before the existing row 376, shifting rows
# 376:415 down to 378:417 and growing the sheet dimension to A1:R417.
$ws.Rows("376:377").Insert()

# Clone the row layout/formatting (including the Primera/Segunda
# labels, units, origin text, etc.) from the prior weekly pair
# (rows 374:375) into the newly inserted rows.
$ws.Range("A374:R375").Copy()
$ws.Range("A376").PasteSpecial()

# Row 376 (Primera) - new weekly price observation
$ws.Range("D376").Value = 44918
$ws.Range("J376").Value = 1200
$ws.Range("K376").Value = 400
$ws.Range("L376").Value = 500
$ws.Range("M376").Value = 442
$ws.Range("P376").Value = 110

# Row 377 (Segunda) - new weekly price observation
$ws.Range("D377").Value = 44918
$ws.Range("J377").Value = 1400
$ws.Range("K377").Value = 400
$ws.Range("L377").Value = 500
$ws.Range("M377").Value = 443
$ws.Range("P377").Value = 89
